$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '57.275.02'
$ws.Range('E2').Value = '  -2.64%  '
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '2.537.48'
$ws.Range('E3').Value = '  -4.32%  '
$ws.Range('E4').Value = '  -0.02%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '513.56'
$ws.Range('E5').Value = '  -1.74%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '139.67'
$ws.Range('E6').Value = '  -3.12%  '
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.999'
$ws.Range('E7').Value = '  -0.01%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.556'
$ws.Range('E8').Value = '  -2.69%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '6.47'
$ws.Range('E9').Value = '  -7.84%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '0.0990'
$ws.Range('E10').Value = '  -3.36%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.324'
$ws.Range('E11').Value = '  -3.01%  '
$ws.Range('E12').Value = '  +0.09%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '2.984.24'
$ws.Range('E13').Value = '  -4.28%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '57.265.40'
$ws.Range('E14').Value = '  -2.71%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '19.96'
$ws.Range('E15').Value = '  -5.08%  '
$ws.Range('E16').Value = '  -2.75%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '2.564.88'
$ws.Range('E17').Value = '  -4.77%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '331.17'
$ws.Range('E18').Value = '  -2.26%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '4.26'
$ws.Range('E19').Value = '  -2.22%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '10.06'
$ws.Range('E20').Value = '  -2.85%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '6.10'
$ws.Range('E21').Value = '  -3.82%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '1.00'
$ws.Range('E22').Value = '  -0.06%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '64.40'
$ws.Range('E23').Value = '  +0.98%  '
$ws.Range('E24').Value = '  -0.30%  '
$ws.Range('E25').Value = '  -0.01%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '0.399'
$ws.Range('E26').Value = '  -4.27%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '2.661.41'
$ws.Range('E27').Value = '  -4.00%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '6.92'
$ws.Range('E28').Value = '  -1.91%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '0.0₃0747'
$ws.Range('E29').Value = '  -6.57%  '
$ws.Range('E30').Value = '  +0.00%  '
$ws.Range('E31').Value = '  -6.57%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '1.55'
$ws.Range('E32').Value = '  -2.28%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '148.77'
$ws.Range('E33').Value = '  -0.44%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '18.38'
$ws.Range('E34').Value = '  -2.32%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '3.95'
$ws.Range('E35').Value = '  -4.47%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '1.12'
$ws.Range('E36').Value = '  -5.41%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '0.839'
$ws.Range('E37').Value = '  -5.65%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '35.71'
$ws.Range('E38').Value = '  -2.78%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '0.822'
$ws.Range('E39').Value = '  -5.28%  '
$ws.Range('E40').Value = '  -4.57%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '1.00'
$ws.Range('E41').Value = '  +0.07%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '3.45'
$ws.Range('E42').Value = '  -3.44%  '
$ws.Range('B43').Value = 'Stellar'
$ws.Range('C43').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '0.0950'
$ws.Range('E43').Value = '  -1.89%  '
$ws.Range('B44').Value = 'WhiteBITCoin'
$ws.Range('C44').Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '10.60'
$ws.Range('E44').Value = '  -0.55%  '
$ws.Range('E45').Value = '  -6.84%  '
$ws.Range('B46').Value = 'Bittensor'
$ws.Range('C46').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '257.64'
$ws.Range('E46').Value = '  -6.38%  '
$ws.Range('B47').Value = 'Hedera'
$ws.Range('C47').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '0.0518'
$ws.Range('E47').Value = '  -2.60%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '18.39'
$ws.Range('E48').Value = '  -7.49%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '1.963.39'
$ws.Range('E49').Value = '  -3.62%  '
$ws.Range('E50').Value = '  -2.90%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '4.48'
$ws.Range('E51').Value = '  -4.98%  '
